{"js": "// Removed placeholder texts, added borders to boxes\n//\n// 1. The \"For CSS\" paragraph's second run (\"and also validated the CSS on \")\n//    is rewritten to describe validating on VSCode Editor plus\n//    atatus.com/tools/css-lint, including the extra commentary about\n//    contrast/redundant headers.\n// 2. A brand new paragraph about the Cross Browser check (sauce labs) is\n//    inserted between the two blank paragraphs that follow the CSS\n//    paragraph.\n\nconst body = context.document.body;\n\n// --- Step 1: rewrite the CSS-validation sentence -------------------------\nconst oldCssRunResults = body.search(\"and also validated the CSS on \", { matchCase: true });\noldCssRunResults.load(\"items\");\nawait context.sync();\n\nif (oldCssRunResults.items.length === 0) {\n  throw new Error(\"Could not find the CSS validation run to update\");\n}\n\nconst newCssText =\n  \"on VSCode Editor \" +\n  \"and validated the CSS on \" +\n  \"atatus.com/tools/css-lint.  There was nothing major found in either one.  \" +\n  \"However, I was checking out the readability of the page and found the \" +\n  \"contrast between some of the colours were off, so I improved that, in \" +\n  \"addition, I removed some other redundant headers that were in that code \" +\n  \"that did not do anything. \";\n\noldCssRunResults.items[0].insertText(newCssText, \"Replace\");\nawait context.sync();\n\n// --- Step 2: insert the new Cross Browser paragraph -----------------------\n// Locate the \"For CSS\" paragraph, then its immediately-following (first)\n// blank paragraph; the new paragraph is inserted right after that blank one,\n// i.e. between the two originally-adjacent blank paragraphs.\nconst cssParaResults = body.search(\"For CSS: I used the Stylelint\", { matchCase: true });\ncssParaResults.load(\"items\");\nawait context.sync();\n\nif (cssParaResults.items.length === 0) {\n  throw new Error(\"Could not find the 'For CSS' paragraph\");\n}\n\nconst cssParagraph = cssParaResults.items[0].paragraphs.getFirst();\nconst firstBlankParagraph = cssParagraph.getNext();\nawait context.sync();\n\nconst crossBrowserText =\n  \"For the Cross Browser check: I used sauce labs, I did not find any \" +\n  \"noticeable difference in functionality or appearance between the \" +\n  \"browsers.  The only thing that I noticed was an animation that \" +\n  \"functioned in Safari, Firefox and Edge, but strangely enough not \" +\n  \"Chrome.  \" +\n  \"I decided to remove it to be consistent across the browsers.\" +\n  \" \";\n\nfirstBlankParagraph.insertParagraph(crossBrowserText, \"After\");\nawait context.sync();\n", "ps1": "# Removed placeholder texts, added borders to boxes\n#\n# 1. The \"For CSS\" paragraph's second run (\"and also validated the CSS on \")\n#    is rewritten to describe validating on VSCode Editor plus\n#    atatus.com/tools/css-lint, including the extra commentary about\n#    contrast/redundant headers.\n# 2. A brand new paragraph about the Cross Browser check (sauce labs) is\n#    inserted between the two blank paragraphs that follow the CSS\n#    paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: rewrite the CSS-validation sentence --------------------------\n$newCssText = \"on VSCode Editor and validated the CSS on atatus.com/tools/css-lint.  There was nothing major found in either one.  However, I was checking out the readability of the page and found the contrast between some of the colours were off, so I improved that, in addition, I removed some other redundant headers that were in that code that did not do anything. \"\n\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\n    \"and also validated the CSS on \",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $newCssText,\n    2\n)\nif (-not $found) {\n    throw \"Could not find the CSS validation run to update\"\n}\n\n# --- Step 2: insert the new Cross Browser paragraph ------------------------\n# Locate the \"For CSS\" paragraph, then its immediately-following (first)\n# blank paragraph; the new paragraph is inserted right after that blank one,\n# i.e. between the two originally-adjacent blank paragraphs.\n$cssParaIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"For CSS: I used the Stylelint\")) {\n        $cssParaIndex = $i\n        break\n    }\n}\nif ($cssParaIndex -eq 0) {\n    throw \"Could not find the 'For CSS' paragraph\"\n}\n\n$crossBrowserText = \"For the Cross Browser check: I used sauce labs, I did not find any noticeable difference in functionality or appearance between the browsers.  The only thing that I noticed was an animation that functioned in Safari, Firefox and Edge, but strangely enough not Chrome.  I decided to remove it to be consistent across the browsers. \"\n\n$firstBlankParaIndex = $cssParaIndex + 1\n$firstBlankPara = $d.Paragraphs.Item($firstBlankParaIndex)\n$firstBlankPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($firstBlankParaIndex + 1)\n$newPara.Range.Text = $crossBrowserText\n"}
